$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename header labels: "<Name>_old" -> "<Name>_FV2310", "<Name>_new" -> "<Name>_FV2404"
#    (the "diff" column header is left untouched)
# ---------------------------------------------------------------------------
$renameMap = @{
  "Segmentname_old"          = "Segmentname_FV2310";
  "Segmentgruppe_old"        = "Segmentgruppe_FV2310";
  "Segment_old"              = "Segment_FV2310";
  "Datenelement_old"         = "Datenelement_FV2310";
  "Segment ID_old"           = "Segment ID_FV2310";
  "Code_old"                 = "Code_FV2310";
  "Qualifier_old"            = "Qualifier_FV2310";
  "Beschreibung_old"         = "Beschreibung_FV2310";
  "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2310";
  "Bedingung_old"            = "Bedingung_FV2310";
  "Segmentname_new"          = "Segmentname_FV2404";
  "Segmentgruppe_new"        = "Segmentgruppe_FV2404";
  "Segment_new"              = "Segment_FV2404";
  "Datenelement_new"         = "Datenelement_FV2404";
  "Segment ID_new"           = "Segment ID_FV2404";
  "Code_new"                 = "Code_FV2404";
  "Qualifier_new"            = "Qualifier_FV2404";
  "Beschreibung_new"         = "Beschreibung_FV2404";
  "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2404";
  "Bedingung_new"            = "Bedingung_FV2404";
}

$headerRange = $ws.Range("A1:U1")
$colCount = $headerRange.Columns.Count
for ($col = 1; $col -le $colCount; $col++) {
  $cell = $ws.Cells.Item(1, $col)
  $cur = $cell.Value2
  if ($renameMap.ContainsKey($cur)) {
    $cell.Value = $renameMap[$cur]
  }
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into an Excel Table ("Table1") with autofilter and
#    banded rows, matching the sheet's used range A1:U67
# ---------------------------------------------------------------------------
$usedRange = $ws.Range("A1:U67")
$lo = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$lo.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1, frozen pane)
# ---------------------------------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
